# Applies a cyclic upward shift to rows 8..19 for the columns that hold the
# per-observation data (A,B,C,D,E,F,G,H,P,Q,R,S,AW,AX). The content that was
# in row 8 moves to row 9, row 9's content moves to row 10, ... and row 19's
# original content wraps around into row 8. All other columns (I, T, U, V,
# W, Y, Z, AA, AB, AD, AE, AG, AT, AY) are identical across these rows and
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 8
$lastRow = 19

# Column indices (1-based) for the columns that participate in the shift.
$cols = @(1, 2, 3, 4, 5, 6, 7, 8, 16, 17, 18, 19, 49, 50)

# Remember the values currently in the last row (row 19); they will wrap
# around into the first row (row 8) once everything else has shifted down.
$wrapValues = @{}
foreach ($c in $cols) {
    $wrapValues[$c] = $ws.Cells.Item($lastRow, $c).Value()
}

# Shift rows down from the bottom up so we never overwrite data before it
# has been read: row (r-1)'s old value becomes row r's new value.
for ($r = $lastRow; $r -gt $firstRow; $r--) {
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($r - 1, $c).Value()
    }
}

# Finally, write the saved last-row values into the first row.
foreach ($c in $cols) {
    $ws.Cells.Item($firstRow, $c).Value = $wrapValues[$c]
}
